$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{Addr='D2'; Val='332.55'},
    @{Addr='E2'; Val='1.39%'},
    @{Addr='D3'; Val='45.89'},
    @{Addr='E3'; Val='4.46%'},
    @{Addr='D4'; Val='5.680'},
    @{Addr='E4'; Val='3.29%'},
    @{Addr='D5'; Val='0.08383'},
    @{Addr='E5'; Val='4.58%'},
    @{Addr='D6'; Val='2.043'},
    @{Addr='E6'; Val='1.53%'},
    @{Addr='B7'; Val='GateToken'},
    @{Addr='C7'; Val='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'},
    @{Addr='D7'; Val='4.483'},
    @{Addr='E7'; Val='3.85%'},
    @{Addr='B8'; Val='MXToken'},
    @{Addr='C8'; Val='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'},
    @{Addr='D8'; Val='0.9917'},
    @{Addr='E8'; Val='4.22%'},
    @{Addr='B9'; Val='BTSEToken'},
    @{Addr='C9'; Val='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'},
    @{Addr='D9'; Val='2.575'},
    @{Addr='E9'; Val='-0.34%'},
    @{Addr='B10'; Val='LiechtensteinCryptoassetsExchange'},
    @{Addr='C10'; Val='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'},
    @{Addr='D10'; Val='0.1155'},
    @{Addr='E10'; Val='2.68%'},
    @{Addr='B11'; Val='WazirX'},
    @{Addr='C11'; Val='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'},
    @{Addr='D11'; Val='0.1931'},
    @{Addr='E11'; Val='3.40%'},
    @{Addr='B12'; Val='MCDex'},
    @{Addr='C12'; Val='https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'},
    @{Addr='D12'; Val='10.41'},
    @{Addr='E12'; Val='-2.05%'},
    @{Addr='B13'; Val='MandalaExchangeToken'},
    @{Addr='C13'; Val='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'},
    @{Addr='D13'; Val='0.09952'},
    @{Addr='E13'; Val='1.47%'},
    @{Addr='B14'; Val='BitrueCoin'},
    @{Addr='C14'; Val='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'},
    @{Addr='D14'; Val='0.04676'},
    @{Addr='E14'; Val='2.46%'},
    @{Addr='B15'; Val='BitMartToken'},
    @{Addr='C15'; Val='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'},
    @{Addr='D15'; Val='0.1058'},
    @{Addr='E15'; Val='-0.81%'},
    @{Addr='B16'; Val='BitForexToken'},
    @{Addr='C16'; Val='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'},
    @{Addr='D16'; Val='0.001282'},
    @{Addr='E16'; Val='0.35%'},
    @{Addr='B17'; Val='TigerCash'},
    @{Addr='C17'; Val='https://coinranking.com/coin/6hIn06L2+tigercash-tch'},
    @{Addr='D17'; Val='0.006111'},
    @{Addr='E17'; Val='3.65%'},
    @{Addr='B18'; Val='LEO'},
    @{Addr='C18'; Val='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'},
    @{Addr='D18'; Val='3.372'},
    @{Addr='E18'; Val='0.50%'},
    @{Addr='E19'; Val='-3.15%'},
    @{Addr='E20'; Val='-0.34%'},
    @{Addr='D21'; Val='0.2653'},
    @{Addr='E21'; Val='4.35%'},
    @{Addr='D22'; Val='0.04205'},
    @{Addr='E22'; Val='3.43%'},
    @{Addr='D23'; Val='0.001311'},
    @{Addr='E23'; Val='4.22%'},
    @{Addr='D24'; Val='0.004631'},
    @{Addr='E24'; Val='6.99%'},
    @{Addr='D25'; Val='0.0001283'},
    @{Addr='E25'; Val='10.73%'},
    @{Addr='D26'; Val='0.0003749'},
    @{Addr='E26'; Val='0.25%'},
    @{Addr='D38'; Val='0.02801'},
    @{Addr='E38'; Val='9.38%'},
    @{Addr='E39'; Val='2.23%'},
    @{Addr='D40'; Val='0.007741'},
    @{Addr='E40'; Val='2.81%'},
    @{Addr='D41'; Val='0.1435'},
    @{Addr='E41'; Val='2.74%'},
    @{Addr='D42'; Val='0.007283'},
    @{Addr='E42'; Val='-4.38%'},
    @{Addr='D43'; Val='0.002112'},
    @{Addr='E43'; Val='4.93%'},
    @{Addr='D44'; Val='0.009051'},
    @{Addr='E44'; Val='2.15%'},
    @{Addr='D45'; Val='0.3406'},
    @{Addr='D46'; Val='0.00007380'},
    @{Addr='E46'; Val='3.95%'},
    @{Addr='D47'; Val='0.00000000752'},
    @{Addr='E47'; Val='0.36%'},
    @{Addr='D48'; Val='0.0005814'},
    @{Addr='E48'; Val='0.04%'},
    @{Addr='B49'; Val='CoinbaseStockToken'},
    @{Addr='C49'; Val='https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'},
    @{Addr='D49'; Val='0.003506'},
    @{Addr='E49'; Val='-0.65%'},
    @{Addr='B50'; Val='BOLO'},
    @{Addr='C50'; Val='https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'},
    @{Addr='D50'; Val='0.003506'},
    @{Addr='E50'; Val='11.89%'},
    @{Addr='D51'; Val='0.00002105'},
    @{Addr='E51'; Val='0.36%'}
)

foreach ($e in $edits) {
    $rng = $ws.Range($e.Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $e.Val
    $rng.ClearFormats()
}
